# Insert a new weekly price record for "Haba" (Macroferia Regional de Talca)
# as row 71, pushing the existing rows 71-93 down to 72-94 (dimension
# A1:R93 -> A1:R94).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 71..93 down by one row, creating a blank row 71.
$ws.Rows.Item(71).Insert()

# Populate the newly inserted row 71 with the new record.
$ws.Cells.Item(71, 1).Value  = 5
$ws.Cells.Item(71, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(71, 3).Value  = "Maule"
$ws.Cells.Item(71, 4).Value  = 44839
$ws.Cells.Item(71, 5).Value  = 7
$ws.Cells.Item(71, 6).Value  = 100112026
$ws.Cells.Item(71, 7).Value  = "Haba"
$ws.Cells.Item(71, 8).Value  = "Sin especificar"
$ws.Cells.Item(71, 9).Value  = "Primera"
$ws.Cells.Item(71, 10).Value = 200
$ws.Cells.Item(71, 11).Value = 9000
$ws.Cells.Item(71, 12).Value = 9000
$ws.Cells.Item(71, 13).Value = 9000
$ws.Cells.Item(71, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(71, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(71, 16).Value = 360
$ws.Cells.Item(71, 17).Value = 25
$ws.Cells.Item(71, 18).Value = "Hortaliza"
